# Update marksheet correction/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: Right column total marking value
$ws.Range("B11").Value = 5

# "Total" row: Right column total score, and Max column Corr/total text
$ws.Range("B12").Value = 45
$ws.Range("E12").Value = "45/140"
